# Resize/restyle the two arrow connectors that link the game-title
# textbox area (connector 4 and connector 16), per the commit diff:
#   - both connectors grow from ~195209 EMU tall to 396000 EMU tall
#   - both connectors get a heavier outline (28575 -> 63500 EMU = 2.25pt -> 5pt)
#   - connector 4 moves up (new off.y = 1685925 EMU)
#   - connector 16 keeps its (flipped) anchor essentially in place

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 9: "直接箭头连接符 4" (straight arrow connector 4)
$conn4 = $s.Shapes.Item(9)
$conn4.Top = 1685925 / 12700
$conn4.Height = 396000 / 12700
$conn4.Line.Weight = 63500 / 12700

# Shape 10: "直接箭头连接符 16" (straight arrow connector 16, flipped vertically)
$conn16 = $s.Shapes.Item(10)
$conn16.Height = 396000 / 12700
$conn16.Top = 173.7978
$conn16.Line.Weight = 63500 / 12700
